$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Video Recording Time = " note to include the recorded time.
$ws.Range("G1").Value = "Video Recording Time = 17:53"

# Row 12's A-cell was a time value; it is now a text note to check the video.
$ws.Range("A12").Value = "Check video"

# Row 13's A-cell previously held a specific time-check note; simplify it to
# the same "Check video" note used elsewhere.
$ws.Range("A13").Value = "Check video"

# Reposition the window's scroll/view position slightly.
$win = $excel.ActiveWindow
$win.Top = 3080
$win.Left = 6320

# Update the active selection on the sheet.
$ws.Range("A7").Select()
